# Renumber the "STORY n" title slides.
# Physical slide order (Slides.Item(N)) -> old title -> new title:
#   Slide 5  : STORY 1  -> STORY 6
#   Slide 6  : STORY 2  -> STORY 7
#   Slide 7  : STORY 3  -> STORY 8
#   Slide 8  : STORY 4  -> STORY 9
#   Slide 9  : STORY 5  -> STORY 10
#   Slide 10 : STORY 6  -> STORY 11
#   Slide 11 : STORY 7  -> STORY 12
#
# Each "STORY n" text lives in Shapes.Item(1) (the Title placeholder) of the
# respective slide. Update slides 10 and 11 first so we never overwrite a
# value we still need to read later.

$p = $ppt.ActivePresentation

$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 11"
$p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 12"

$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 6"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 7"
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 8"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 9"
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Text = "STORY 10"
